# Commit: "Consistently use double-quotes for all JSON snippets"
#
# The "relayout argument" code sample textbox on slide 1 mixes single and
# double quotes in its two JSON-ish snippets:
#   {'xaxis.range': [0, 1],
#    'yaxis.range': [1.5, 2.5]}
# Both should use double quotes:
#   {"xaxis.range": [0, 1],
#    "yaxis.range": [1.5, 2.5]}
#
# We edit the existing runs in place (same run/paragraph boundaries, same
# rPr) so only the four single-quote characters turn into double-quotes -
# exactly what the canonical-XML diff shows.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the textbox holding the "relayout argument:" code sample.
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $t = $shp.TextFrame.TextRange.Text
        if ($t.Contains("relayout") -and $t.Contains("xaxis.range")) {
            $targetShape = $shp
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Each entry is the 1-based start offset (counting paragraph marks as one
# character each, matching TextRange/Characters indexing), the run length,
# the expected existing text (sanity check) and the corrected text.
$edits = @(
    @(20, 2,  "{'",             '{"'),
    @(33, 10, "': [0, 1],",     '": [0, 1],'),
    @(44, 2,  " '",             ' "'),
    @(57, 14, "': [1.5, 2.5]}", '": [1.5, 2.5]}')
)

foreach ($edit in $edits) {
    $start = $edit[0]
    $len = $edit[1]
    $expected = $edit[2]
    $replacement = $edit[3]

    $run = $tr.Characters($start, $len)
    if ($run.Text -ne $expected) {
        Write-Host ("Unexpected text at " + $start + "," + $len + ": [" + $run.Text + "] expected [" + $expected + "]")
    }
    $run.Text = $replacement
}

Write-Host ("Final text: [" + $tr.Text + "]")
